$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date formatting from an existing date cell (B2) down to the
# new date cells (B9:B13) so they reuse the existing date style.
$ws.Range("B2").Copy()
$ws.Range("B9:B13").PasteSpecial(-4122)

# New rows of data (rows 9-13), continuing the existing ID/date/test pattern
$data = @(
    @(8,  43970, 1, 0, 1, 0, 0),
    @(9,  43971, 0, 1, 0, 0, 0),
    @(10, 43972, 0, 0, 0, 0, 1),
    @(11, 43973, 0, 0, 0, 1, 0),
    @(12, 43974, 0, 0, 1, 0, 0)
)

$r = 9
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $r++
}

$ws.Range("A14").Select()
